# Add info about two new JK (housing complexes): "Дом Соболева" and "SOKOLNIKI"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 6: Дом Соболева ----
$ws.Range("A6").Value = "Дом Соболева"
$ws.Range("B6").Value = 33
$ws.Range("C6").Value = 0
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 9
$ws.Range("F6").Value = 12
$ws.Range("G6").Value = 12
$ws.Range("H6").Value = 2
$ws.Range("I6").Value = 10
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("M6").Formula = "=(J6+K6+L6)/B6"
$ws.Range("N6").Formula = "=I6-J6"
$ws.Range("O6").Formula = "=33 * 146"
$ws.Range("P6").Value = 146
$ws.Range("Q6").Value = 3.5
$ws.Range("R6").Value = 3.95
$ws.Range("S6").Value = 10
$ws.Range("T6").Value = 10
$ws.Range("U6").Value = 6
$ws.Range("V6").Value = 2
$ws.Range("W6").Formula = "=U6/V6"
$ws.Range("X6").Value = 0
$ws.Range("Y6").Value = 0
$ws.Range("Z6").Value = 0
$ws.Range("AA6").Value = 1
$ws.Range("AB6").Value = 1
$ws.Range("AC6").Value = 0
$ws.Range("AD6").Value = 0
$ws.Range("AE6").Value = 2
$ws.Range("AF6").Value = 55.764016
$ws.Range("AG6").Value = 37.652657

# ---- Row 7: SOKOLNIKI ----
$ws.Range("A7").Value = "SOKOLNIKI"
$ws.Range("B7").Value = 763
$ws.Range("C7").Value = 0
$ws.Range("D7").Value = 325
$ws.Range("E7").Value = 229
$ws.Range("F7").Value = 181
$ws.Range("G7").Value = 28
$ws.Range("H7").Value = 9
$ws.Range("I7").Value = 986
$ws.Range("J7").Value = 313
$ws.Range("K7").Value = 85
$ws.Range("L7").Value = 0
$ws.Range("M7").Formula = "=(J7+K7+L7)/B7"
$ws.Range("N7").Formula = "=I7-J7"
$ws.Range("O7").Value = 45984
$ws.Range("P7").Value = 60.3
$ws.Range("Q7").Value = 3
$ws.Range("R7").Value = 3
$ws.Range("S7").Value = 36
$ws.Range("T7").Value = 36
$ws.Range("U7").Value = 16
$ws.Range("V7").Value = 3
$ws.Range("W7").Formula = "=U7/V7"
$ws.Range("X7").Value = 3
$ws.Range("Y7").Value = 2
$ws.Range("Z7").Value = 0
$ws.Range("AA7").Value = 1
$ws.Range("AB7").Value = 2
$ws.Range("AC7").Value = 1
$ws.Range("AD7").Value = 0
$ws.Range("AE7").Value = 0
$ws.Range("AF7").Value = 55.782995
$ws.Range("AG7").Value = 37.690437

# Update view state
$ws.Range("AG7").Select()
